# Make the letter's static placeholder-ish text dynamic by substituting
# template tokens, per the commit "feat: dynamis content latter".

$d = $word.ActiveDocument

# 1) "No  :    YDDS/II/D/{no_surat}" -> "No  :    {no_surat}"
$d.Content.Find.Execute(":    YDDS/II/D/{no_surat}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ":    {no_surat}", 2)

# 2) "Jakarta, {date}" -> "{latter_province}, {date}"
$d.Content.Find.Execute("Jakarta, {date}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{latter_province}, {date}", 2)

# 3) "Kepada yth. :" -> "{greeting}"
$d.Content.Find.Execute("Kepada yth. :", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{greeting}", 2)

# 4) "Surat Saudara Tgl……….​" -> "Surat Saudara Tgl {date}"
$d.Content.Find.Execute("Surat Saudara Tgl" + [char]8230 + [char]8230 + [char]8230 + ".", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Surat Saudara Tgl {date}", 2)

# 5) "... pensiunan a/n. Saudara." -> "... pensiunan a/n. {name}."
$d.Content.Find.Execute("pensiunan a/n. Saudara.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "pensiunan a/n. {name}.", 2)

# 6) "Pensiunan a/n. Saudara tsb." -> "Pensiunan a/n. {name} tsb."
$d.Content.Find.Execute("Pensiunan a/n. Saudara tsb.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pensiunan a/n. {name} tsb.", 2)

# 7) "belum dilengkapi dengan kuitansi asli dari Rumah Sakit." -> "{short_desc}."
$d.Content.Find.Execute("belum dilengkapi dengan kuitansi asli dari Rumah Sakit.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{short_desc}.", 2)

# 8) Long paragraph -> "Berdasarkan dengan hal-hal tersebut atas, Kami harapkan agar {long_desc}."
$d.Content.Find.Execute("Berdasarkan dengan hal-hal tersebut atas, Kami harapkan agar dilengkapi dengan kuitansi asli dari Rumah Sakit, klaim ulang (untuk melengkapi berkas) diterima YDDS selambat - lambatnya 60 hari sejak tanggal pemberitahuan/surat YDDS.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Berdasarkan dengan hal-hal tersebut atas, Kami harapkan agar {long_desc}.", 2)

# The paragraph immediately following the long-description paragraph was an
# empty spacer paragraph; the edit removes it entirely (merging it away).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Berdasarkan dengan hal-hal tersebut atas, Kami harapkan agar {long_desc}.")) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
